# Update countries & provincias Spain
# Refresh country statistics and re-rank rows that changed order, plus
# update the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 12:35"

# Row 27 - numeric refresh only (country unchanged)
$ws.Range("B27").Value = 30618
$ws.Range("C27").Value = 21
$ws.Range("E27").Value = 1132

# Rows 32/33 - Emiratos Arabes Unidos overtakes Irlanda
$ws.Range("A32").Value = "Emiratos Arabes Unidos"
$ws.Range("B32").Value = 25063
$ws.Range("C32").Value = 873
$ws.Range("D32").Value = 10791
$ws.Range("E32").Value = 14045
$ws.Range("G32").Value = 3
$ws.Range("H32").Value = 227

$ws.Range("A33").Value = "Irlanda"
$ws.Range("B33").Value = 24200
$ws.Range("D33").Value = 19470
$ws.Range("E33").Value = 3183
$ws.Range("H33").Value = 1547

# Row 37 - numeric refresh only (country unchanged)
$ws.Range("D37").Value = 10166
$ws.Range("E37").Value = 5899
$ws.Range("G37").Value = 6
$ws.Range("H37").Value = 1126

# Rows 55/56 - Barein overtakes Argelia
$ws.Range("A55").Value = "Barein"
$ws.Range("B55").Value = 7374
$ws.Range("C55").Value = 190
$ws.Range("D55").Value = 2952
$ws.Range("E55").Value = 4410
$ws.Range("H55").Value = 12

$ws.Range("A56").Value = "Argelia"
$ws.Range("B56").Value = 7201
$ws.Range("D56").Value = 3625
$ws.Range("E56").Value = 3021
$ws.Range("H56").Value = 555

# Row 59 - numeric refresh only (country unchanged)
$ws.Range("B59").Value = 6972
$ws.Range("C59").Value = 20
$ws.Range("D59").Value = 3890
$ws.Range("E59").Value = 2889
$ws.Range("G59").Value = 1
$ws.Range("H59").Value = 193

# Row 96 - numeric refresh only (country unchanged)
$ws.Range("B96").Value = 1467
$ws.Range("C96").Value = 1
$ws.Range("E96").Value = 28

# Row 107 - numeric refresh only (country unchanged)
$ws.Range("B107").Value = 949
$ws.Range("C107").Value = 1
$ws.Range("D107").Value = 742
$ws.Range("E107").Value = 176

# Row 136 - numeric refresh only (country unchanged)
$ws.Range("B136").Value = 365
$ws.Range("C136").Value = 13
$ws.Range("D136").Value = 120
$ws.Range("E136").Value = 240

# Rows 143/144/145 - Madagascar overtakes Vietnam and Montenegro
$ws.Range("A143").Value = "Madagascar"
$ws.Range("B143").Value = 326
$ws.Range("C143").Value = 4
$ws.Range("D143").Value = 119
$ws.Range("E143").Value = 205
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 2

$ws.Range("A144").Value = "Vietnam"
$ws.Range("B144").Value = 324
$ws.Range("D144").Value = 263
$ws.Range("E144").Value = 61
$ws.Range("H144").Value = 0

$ws.Range("A145").Value = "Montenegro"
$ws.Range("B145").Value = 324
$ws.Range("D145").Value = 311
$ws.Range("E145").Value = 4
$ws.Range("H145").Value = 9
